# "did FJLT 2 for d = 64"
# Adds a new Sheet2 after Sheet1 with a small scratch table computing
# d^(1/3) and SQRT(d) for a handful of d values, formatted to 2 decimals.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 and make it the active sheet.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Headers (bold, matching the style used for Sheet1's headers)
$ws2.Range("B1").Value = "FJLT 1"
$ws2.Range("C1").Value = "FJLT 2"
$ws2.Range("B1:C1").Font.Bold = $true

$ws2.Range("A2").Value = "d"
$ws2.Range("B2").Value = "d^(1/3)"
$ws2.Range("C2").Value = "d^(1/2)"
$ws2.Range("A2:C2").Font.Bold = $true

# Data rows: d = 32, 64, 128, 256, 512, 1024
$dValues = 32, 64, 128, 256, 512, 1024
$row = 3
foreach ($d in $dValues) {
    $ws2.Cells.Item($row, 1).Value = $d
    $ws2.Cells.Item($row, 2).Formula = "=A$row^(1/3)"
    $ws2.Cells.Item($row, 3).Formula = "=SQRT(A$row)"
    $row++
}

# Format the formula results (B3:C8) to 2 decimal places.
$ws2.Range("B3:C8").NumberFormat = "0.00"

$ws2.Range("F10").Select()

$wb.Save()
